$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$tr = $s.Shapes.Item(2).TextFrame.TextRange

# --- Paragraph 1: split "...creator of LaTeX. The Part-Time Parliament"
# into "...creator of LaTeX" + ". "
$para1 = $tr.Paragraphs(1)
$found = $para1.Find("LaTeX")
$tailStart = $found.Start + $found.Length
$tailLen = ($para1.Start + $para1.Length) - $tailStart
$tail = $tr.Characters($tailStart, $tailLen)
$tail.Text = ". "

# --- Paragraph 2: split "The paper title was ...publish in ACM 1998."
# into "The " + "paper title was ...publish in ACM 1998."
$para2 = $tr.Paragraphs(2)
$lead = $para2.Characters(1, 4)
$lead.Text = "The "
